$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places
$row5Values = @{
    "B5" = 12.01; "C5" = 8.68; "D5" = 0.96; "E5" = 26.08; "F5" = 21.08; "G5" = 9.45;
    "H5" = 33.22; "I5" = 14.54; "J5" = 6.34; "K5" = 9.34; "L5" = 10.47; "M5" = 10.99;
    "N5" = 3.02; "O5" = 9.4; "P5" = 13.25; "Q5" = 8.06; "R5" = 0.77; "S5" = 0.59;
    "T5" = 135.7; "U5" = 26.1; "V5" = 8.68; "W5" = 17.37; "X5" = 9.14; "Y5" = 1.56;
    "Z5" = 16.26; "AA5" = 7.66; "AB5" = 6.88; "AC5" = 8.07; "AD5" = 10.95; "AE5" = 0.55;
    "AF5" = 29.77; "AG5" = 4.81; "AH5" = 10.85
}

foreach ($addr in $row5Values.Keys) {
    $ws.Range($addr).Value = $row5Values[$addr]
}

# Delete row 6 entirely
$ws.Rows.Item(6).Delete()
